# Update the "sun profile" (grid connection) yearly values in column C
# and move the active selection on the "Yearly" sheet, per the diff:
#   C2: 446.20530973451326 -> 531
#   C3: 601.64302477183833 -> 782
#   C4: 749.2263709677419  -> 975
#   C5: 755.63702857142857 -> 881
#   C6: 819.83133241953021 -> 1000
#   selection moves from F10 to E6

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Yearly")

$ws.Range("C2").Value = 531
$ws.Range("C3").Value = 782
$ws.Range("C4").Value = 975
$ws.Range("C5").Value = 881
$ws.Range("C6").Value = 1000

$ws.Activate()
[void]$ws.Range("E6").Select()
